# Add Test Results for the web application login test case (sheet "Login").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- Capitalize / correct existing test-case wording (rows 12-18) ---
$ws.Cells.Item(12, 1).Value = "Enter a valid email address and a valid password used in Yammer to login"
$ws.Cells.Item(12, 2).Value = "System should awllow the user to login if they are in HR"
$ws.Cells.Item(12, 3).Value = "System allows user to login"

$ws.Cells.Item(13, 1).Value = "Click the log off button and click to login again"
$ws.Cells.Item(13, 3).Value = "System allows user to login"

$ws.Cells.Item(14, 1).Value = "Enter an invalid email address and a valid password used in Yammer to login"
$ws.Cells.Item(14, 3).Value = "System does not allow the user to login"

$ws.Cells.Item(15, 1).Value = "Enter a valid email address and an invalid password in Yammer to login"
$ws.Cells.Item(15, 3).Value = "System denies login"

$ws.Cells.Item(16, 1).Value = "Enter invalid email address and invalid password to login with Yammer"

$ws.Cells.Item(17, 1).Value = "Press login in Yammer with email and pass word not entered"

$ws.Cells.Item(18, 1).Value = "Press login in the wab app login page "

# --- Row height tweaks for rows 10-12 (explicit custom heights) ---
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 45

# --- New test result rows 19 and 20 ---
$ws.Cells.Item(19, 1).Value = "Bypassing the login page by typing the url of the home page"
$ws.Cells.Item(19, 2).Value = "The system should redirect to the login page"
$ws.Cells.Item(19, 3).Value = "The system redirects to login page"

$ws.Cells.Item(20, 1).Value = "Bypassing the login page to any other page(excludeing home page)."
$ws.Cells.Item(20, 2).Value = "The system should redirect to login page"

$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 30

# Format new rows 19:20 like the table above them (left/right border, wrap, vertical center).
# Only touch cells that actually hold data: A19:C19 and A20:B20 (row 20 has no C cell).
$dataCells = @(
    $ws.Cells.Item(19, 1), $ws.Cells.Item(19, 2), $ws.Cells.Item(19, 3),
    $ws.Cells.Item(20, 1), $ws.Cells.Item(20, 2)
)
foreach ($cell in $dataCells) {
    $cell.WrapText = $true
    $cell.VerticalAlignment = -4108
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(10).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1
}

# C19 gets a light fill to match the "actual result" column styling used elsewhere
$fillSrc = $wb.Worksheets.Item("OpenBody").Range("D17")
$fillSrc.Copy()
$ws.Cells.Item(19, 3).PasteSpecial(-4122)
$ws.Cells.Item(19, 3).WrapText = $true
$ws.Cells.Item(19, 3).VerticalAlignment = -4108
$ws.Cells.Item(19, 3).Borders.Item(8).LineStyle = -4142
$ws.Cells.Item(19, 3).Borders.Item(9).LineStyle = -4142
$ws.Cells.Item(19, 3).Borders.Item(7).Weight = 2
$ws.Cells.Item(19, 3).Borders.Item(7).LineStyle = 1
$ws.Cells.Item(19, 3).Borders.Item(10).Weight = 2
$ws.Cells.Item(19, 3).Borders.Item(10).LineStyle = 1

# --- sheet view: scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L18").Select()
